# IFRS company_list sheet — refresh row 2-9 financial figures.
# (commit: "error solve ifrs list")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("company_list")

# ---- Row 2 : 2014/12 (IFRS연결) ----
$ws.Range("D2").Value  = 9852
$ws.Range("E2").Value  = 125
$ws.Range("F2").Value  = 125
$ws.Range("G2").Value  = 135
$ws.Range("H2").Value  = 88
$ws.Range("I2").Value  = 95
$ws.Range("J2").Value  = -7
$ws.Range("K2").Value  = 68314
$ws.Range("L2").Value  = 60618
$ws.Range("M2").Value  = 7696
$ws.Range("N2").Value  = 7676
$ws.Range("O2").Value  = 20
$ws.Range("P2").Value  = 4408
$ws.Range("Q2").Value  = -2580
$ws.Range("R2").Value  = 2978
$ws.Range("S2").Value  = -1756
$ws.Range("T2").Value  = 28
$ws.Range("U2").ClearContents()
$ws.Range("V2").Value  = 7007
$ws.Range("W2").Value  = 1.26
$ws.Range("X2").Value  = 0.9
$ws.Range("Y2").Value  = 1.25
$ws.Range("Z2").Value  = 0.13
$ws.Range("AA2").Value = 787.66
$ws.Range("AB2").Value = 81.08
$ws.Range("AC2").Value = 93
$ws.Range("AD2").Value = 34.49
$ws.Range("AE2").Value = 7940
$ws.Range("AF2").Value = 0.4
$ws.Range("AG2").Value = 60
$ws.Range("AH2").Value = 1.87
$ws.Range("AI2").Value = 63.32
$ws.Range("AJ2").Value = 97137759

# ---- Row 3 : 2015/12 (IFRS연결) ----
$ws.Range("D3").Value  = 18338
$ws.Range("E3").Value  = -166
$ws.Range("F3").Value  = -166
$ws.Range("G3").Value  = -110
$ws.Range("H3").Value  = -123
$ws.Range("I3").Value  = -122
$ws.Range("J3").Value  = -1
$ws.Range("K3").Value  = 76623
$ws.Range("L3").Value  = 68797
$ws.Range("M3").Value  = 7826
$ws.Range("N3").Value  = 7805
$ws.Range("O3").Value  = 21
$ws.Range("P3").Value  = 4408
$ws.Range("Q3").Value  = 313
$ws.Range("R3").Value  = 1962
$ws.Range("S3").Value  = -764
$ws.Range("T3").Value  = 38
$ws.Range("U3").ClearContents()
$ws.Range("V3").Value  = 5936
$ws.Range("W3").Value  = -0.91
$ws.Range("X3").Value  = -0.67
$ws.Range("Y3").Value  = -1.58
$ws.Range("Z3").Value  = -0.17
$ws.Range("AA3").Value = 879.0599999999999
$ws.Range("AB3").Value = 79.86
$ws.Range("AC3").Value = -120
$ws.Range("AD3").Value = -28.2
$ws.Range("AE3").Value = 7772
$ws.Range("AF3").Value = 0.44
$ws.Range("AG3").Value = 0
$ws.Range("AH3").Value = 0
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 97137759

# ---- Row 4 : 2016/12 (IFRS연결) ----
$ws.Range("D4").Value  = 15788
$ws.Range("E4").Value  = -1929
$ws.Range("F4").Value  = -1929
$ws.Range("G4").Value  = -1860
$ws.Range("H4").Value  = -1608
$ws.Range("I4").Value  = -1609
$ws.Range("J4").Value  = 2
$ws.Range("K4").Value  = 69357
$ws.Range("L4").Value  = 61206
$ws.Range("M4").Value  = 8151
$ws.Range("N4").Value  = 8134
$ws.Range("O4").Value  = 17
$ws.Range("P4").Value  = 8862
$ws.Range("Q4").Value  = -7356
$ws.Range("R4").Value  = 2168
$ws.Range("S4").Value  = 3678
$ws.Range("T4").Value  = 20
$ws.Range("U4").ClearContents()
$ws.Range("V4").Value  = 7068
$ws.Range("W4").Value  = -12.22
$ws.Range("X4").Value  = -10.18
$ws.Range("Y4").Value  = -20.2
$ws.Range("Z4").Value  = -2.21
$ws.Range("AA4").Value = 750.9299999999999
$ws.Range("AB4").Value = -7
$ws.Range("AC4").Value = -1352
$ws.Range("AD4").Value = -1.54
$ws.Range("AE4").Value = 4623
$ws.Range("AF4").Value = 0.45
$ws.Range("AG4").Value = 0
$ws.Range("AH4").Value = 0
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 172442511

# ---- Row 5 : 2017/12 (IFRS연결) ----
$ws.Range("D5").Value  = 17128
$ws.Range("E5").Value  = 655
$ws.Range("F5").Value  = 655
$ws.Range("G5").Value  = 728
$ws.Range("H5").Value  = 557
$ws.Range("I5").Value  = 555
$ws.Range("J5").Value  = 2
$ws.Range("K5").Value  = 69053
$ws.Range("L5").Value  = 60281
$ws.Range("M5").Value  = 8772
$ws.Range("N5").Value  = 8753
$ws.Range("O5").Value  = 18
$ws.Range("P5").Value  = 8862
$ws.Range("Q5").Value  = -3213
$ws.Range("R5").Value  = -141
$ws.Range("S5").Value  = 3477
$ws.Range("T5").Value  = 21
$ws.Range("U5").ClearContents()
$ws.Range("V5").Value  = 11298
$ws.Range("W5").Value  = 3.82
$ws.Range("X5").Value  = 3.25
$ws.Range("Y5").Value  = 6.57
$ws.Range("Z5").Value  = 0.8
$ws.Range("AA5").Value = 687.22
$ws.Range("AB5").Value = 0.01
$ws.Range("AC5").Value = 313
$ws.Range("AD5").Value = 8.619999999999999
$ws.Range("AE5").Value = 4975
$ws.Range("AF5").Value = 0.54
$ws.Range("AG5").Value = 0
$ws.Range("AH5").Value = 0
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 172442511

# ---- Row 6 : 2018/12 (IFRS연결) ----
$ws.Range("D6").Value  = 19019
$ws.Range("E6").Value  = 972
$ws.Range("F6").Value  = 972
$ws.Range("G6").Value  = 974
$ws.Range("H6").Value  = 724
$ws.Range("I6").Value  = 725
$ws.Range("K6").Value  = 73441
$ws.Range("L6").Value  = 63938
$ws.Range("M6").Value  = 9503
$ws.Range("N6").Value  = 9485
$ws.Range("P6").Value  = 8862
$ws.Range("Q6").Value  = -2509
$ws.Range("R6").Value  = -1538
$ws.Range("S6").Value  = 4254
$ws.Range("T6").Value  = 35
$ws.Range("U6").ClearContents()
$ws.Range("V6").Value  = 13240
$ws.Range("W6").Value  = 5.11
$ws.Range("X6").Value  = 3.81
$ws.Range("Y6").Value  = 7.95
$ws.Range("Z6").Value  = 1.02
$ws.Range("AA6").Value = 672.8099999999999
$ws.Range("AB6").Value = 8.27
$ws.Range("AC6").Value = 409
$ws.Range("AD6").Value = 4.98
$ws.Range("AE6").Value = 5391
$ws.Range("AF6").Value = 0.38
$ws.Range("AG6").Value = 0
$ws.Range("AH6").Value = 0
$ws.Range("AI6").ClearContents()
$ws.Range("AJ6").Value = 172442511

# ---- Rows 7-9 : 2019/12(E), 2020/12(E), 2021/12(E) — estimates removed, data wiped ----
$ws.Range("D7:AJ9").ClearContents()
